$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "angekommen" -------------------------------------------------

# Header cell G1, styled like the other header cells (bold font + bottom border),
# matching the look of A1.
$ws.Range("G1").Value = "angekommen"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").Borders.Item(9).LineStyle = 1   # xlEdgeBottom thin border

# Column width for the new column (input value accounts for Excel's internal
# character-width padding so the stored width comes out to 15).
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666

# Copy "Anzahl bestellt" (column F) values into the new "angekommen" column (G)
# for every row that currently has data, and apply the new highlight style
# (green fill, centered horizontally). Alignment is applied before the fill
# color so the style table doesn't end up with a stray intermediate entry.
$rows = @(2,3,4,5,6,7,9,11,12,13,15,16)
foreach ($r in $rows) {
    $srcCell = $ws.Cells.Item($r, 6)
    $dstCell = $ws.Cells.Item($r, 7)
    $dstCell.Value = $srcCell.Value2
    $dstCell.HorizontalAlignment = -4108  # xlCenter
    $dstCell.Interior.Color = 5296274     # RGB(146, 208, 80) == FF92D050
}

# Restore the selection as saved in the workbook
$ws.Range("G7").Select()
